# "Fruta / hortaliza, semanal"
# Insert one new weekly price-report row at row 70 (pushing the existing
# rows 70-94 down to 71-95), and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(70).Insert()

$ws.Range("A70").Value = 8
$ws.Range("B70").Value = "Terminal La Palmera de La Serena"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44726
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112052
$ws.Range("G70").Value = "Albahaca"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 800
$ws.Range("K70").Value = 3500
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = 3750
$ws.Range("N70").Value = "$/paquete"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 3750
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

Write-Output "Inserted new row 70; sheet now spans to row 95."
